$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text, preserving the cells original style/number format.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range("D2") '71.009.68'
$ws.Range("E2").Value = '  +6.01%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.657.23'
$ws.Range("E3").Value = '  +17.89%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
Set-TextValue $ws.Range("D5") '626.99'
$ws.Range("E5").Value = '  +8.67%  '

# Row 6
Set-TextValue $ws.Range("D6") '181.97'
$ws.Range("E6").Value = '  +1.93%  '

# Row 7
Set-TextValue $ws.Range("D7") '3.652.58'
$ws.Range("E7").Value = '  +17.83%  '

# Row 8
$ws.Range("E8").Value = '  -0.08%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.539'
$ws.Range("E9").Value = '  +4.86%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.164'
$ws.Range("E10").Value = '  +8.51%  '

# Row 11
Set-TextValue $ws.Range("D11") '6.70'
$ws.Range("E11").Value = '  +5.68%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.502'
$ws.Range("E12").Value = '  +7.36%  '

# Row 13
Set-TextValue $ws.Range("D13") '40.58'
$ws.Range("E13").Value = '  +12.21%  '

# Row 14
Set-TextValue $ws.Range("D14") '0.0000255'

# Row 15
Set-TextValue $ws.Range("D15") '4.268.19'
$ws.Range("E15").Value = '  +17.90%  '

# Row 16
Set-TextValue $ws.Range("D16") '70.938.63'
$ws.Range("E16").Value = '  +5.97%  '

# Row 17
Set-TextValue $ws.Range("D17") '3.660.28'
$ws.Range("E17").Value = '  +18.04%  '

# Row 18
$ws.Range("E18").Value = '  +1.61%  '

# Row 19
Set-TextValue $ws.Range("D19") '7.56'
$ws.Range("E19").Value = '  +7.73%  '

# Row 20
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D20") '16.94'
$ws.Range("E20").Value = '  +1.58%  '

# Row 21
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range("D21") '519.70'
$ws.Range("E21").Value = '  +7.97%  '

# Row 22
Set-TextValue $ws.Range("D22") '9.27'
$ws.Range("E22").Value = '  +19.82%  '

# Row 23
Set-TextValue $ws.Range("D23") '0.746'
$ws.Range("E23").Value = '  +7.97%  '

# Row 24
Set-TextValue $ws.Range("D24") '88.55'
$ws.Range("E24").Value = '  +5.79%  '

# Row 25
Set-TextValue $ws.Range("D25") '2.50'
$ws.Range("E25").Value = '  +10.87%  '

# Row 26
Set-TextValue $ws.Range("D26") '13.53'
$ws.Range("E26").Value = '  +6.95%  '

# Row 27
Set-TextValue $ws.Range("D27") '11.05'
$ws.Range("E27").Value = '  +9.15%  '

# Row 28
$ws.Range("E28").Value = '  -0.06%  '

# Row 29
Set-TextValue $ws.Range("D29") '2.56'
$ws.Range("E29").Value = '  +12.35%  '

# Row 30
Set-TextValue $ws.Range("D30") '2.94'
$ws.Range("E30").Value = '  +12.89%  '

# Row 31
Set-TextValue $ws.Range("D31") '8.11'
$ws.Range("E31").Value = '  +1.31%  '

# Row 32
Set-TextValue $ws.Range("D32") '31.65'
$ws.Range("E32").Value = '  +12.97%  '

# Row 33
$ws.Range("E33").Value = '  +17.31%  '

# Row 34
$ws.Range("E34").Value = '  +3.55%  '

# Row 35
$ws.Range("E35").Value = '  -0.05%  '

# Row 36
Set-TextValue $ws.Range("D36") '6.14'
$ws.Range("E36").Value = '  +9.68%  '

# Row 37
$ws.Range("E37").Value = '  +8.93%  '

# Row 38
$ws.Range("E38").Value = '  +11.33%  '

# Row 39
Set-TextValue $ws.Range("D39") '2.22'
$ws.Range("E39").Value = '  +10.96%  '

# Row 40
Set-TextValue $ws.Range("D40") '51.94'
$ws.Range("E40").Value = '  +5.82%  '

# Row 41
$ws.Range("E41").Value = '  +5.94%  '

# Row 42
Set-TextValue $ws.Range("D42") '45.86'
$ws.Range("E42").Value = '  -5.47%  '

# Row 43
Set-TextValue $ws.Range("D43") '8.84'
$ws.Range("E43").Value = '  +5.94%  '

# Row 44
Set-TextValue $ws.Range("D44") '3.123.81'
$ws.Range("E44").Value = '  +11.56%  '

# Row 45
Set-TextValue $ws.Range("D45") '425.21'
$ws.Range("E45").Value = '  +13.73%  '

# Row 46
Set-TextValue $ws.Range("D46") '2.78'
$ws.Range("E46").Value = '  +3.64%  '

# Row 47
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D47") '28.67'
$ws.Range("E47").Value = '  +14.11%  '

# Row 48
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D48") '0.0371'
$ws.Range("E48").Value = '  +7.89%  '

# Row 49
Set-TextValue $ws.Range("D49") '139.94'
$ws.Range("E49").Value = '  +3.20%  '

# Row 51
Set-TextValue $ws.Range("D51") '2.47'
$ws.Range("E51").Value = '  +9.73%  '
